$wb = $excel.ActiveWorkbook

# --- CApULAbIFM sheet: replace formula with hardcoded value from U.S. file ---
$wsData = $wb.Worksheets.Item("CApULAbIFM")
$wsData.Range("B2").Value = 1500000

# --- About sheet: update source text, clear now-unused rows ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B3").Value = "Consultation with American Forest Foundation"
$wsAbout.Range("B4").Value = $null
$wsAbout.Range("B5").Value = $null
$wsAbout.Range("B6").Value = $null
$wsAbout.Range("B7").Value = $null
$wsAbout.Hyperlinks.Delete()
$wsAbout.Rows.Item(10).Delete()

# --- Remove the now-unused Calculations sheet ---
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsCalc.Delete()
